# Update the "Förändrad" (changed) date column (C) for rows 2-14
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = [DateTime]::FromOADate(45233)

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
